$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G needs a custom width matching column E (closest width the
# engine's character-width quantization can reach to the source 17.1640625).
$ws.Columns("G").ColumnWidth = 16.33

# Add the new row (19) first so the new shared string "SARIMA_2_xreg" gets
# shared-string index 22, matching the order strings were introduced in the
# target workbook (SARIMA_2_xreg, then Price_MAPE, then Demand_MAPE).
$ws.Range("A19").Value = "SARIMA_2_xreg"

# New header cells for the two new metric columns.
$ws.Range("F1").Value = "Price_MAPE"
$ws.Range("G1").Value = "Demand_MAPE"

# Row 3 (ETS)
$ws.Range("F3").Value = 0.3399347
$ws.Range("G3").Value = 0.04089591

# Row 4 (ETS_log) - F4 keeps the same explicit-black-font style as the
# other already-styled cells in that row (style index 1).
$ws.Range("F4").Value = 1
$ws.Range("F4").Font.Color = 0
$ws.Range("G4").Value = 0.04089592

# Row 5 (SARIMA)
$ws.Range("F5").Value = 0.3273797
$ws.Range("G5").Value = 0.03829565

# Row 7 (SARIMA_log)
$ws.Range("F7").Value = 0.3002373
$ws.Range("G7").Value = 0.03841707

# Row 8 (SARIMA_weather_log) - full new set of metrics including the brand
# new H/I/J/K columns.
$ws.Range("F8").Value = 0.2998991
$ws.Range("G8").Value = 0.03859608
$ws.Range("H8").Value = 12.63273
$ws.Range("I8").Value = 273.0789
$ws.Range("J8").Value = 14438.61
$ws.Range("K8").Value = 330403616

# Row 19 (SARIMA_2_xreg) - new row with all new metrics.
$ws.Range("F19").Value = 0.3268286
$ws.Range("G19").Value = 0.03843157
$ws.Range("H19").Value = 12.49192
$ws.Range("I19").Value = 268.3188
$ws.Range("J19").Value = 14369.06
$ws.Range("K19").Value = 327475674

# Match the author's final cursor/selection position recorded in the file.
$ws.Range("H15").Select() | Out-Null
